# Update "Horas Extra" (H) and "Salario Total" (J) columns on the
# "Datos Empleados" sheet for the employee rows (2-22).
#
# The source file stores every cell as text (inline strings), so the
# replacement values are written back as text as well (using the classic
# leading-apostrophe trick) to keep them looking the same as the rest of
# the column, e.g. "170.0" rather than the number 170.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Empleados")

$updates = @{
    2  = @{ H = "0"; J = "170.0" }
    3  = @{ H = "0"; J = "170.0" }
    4  = @{ H = "0"; J = "2040.0" }
    5  = @{ H = "0"; J = "20400.0" }
    6  = @{ H = "0"; J = "187.425" }
    7  = @{ H = "0"; J = "14960.0" }
    8  = @{ H = "0"; J = "204.0" }
    9  = @{ H = "0"; J = "244.8" }
    10 = @{ H = "0"; J = "244.8" }
    11 = @{ H = "0"; J = "340.0" }
    12 = @{ H = "0"; J = "170000.0" }
    13 = @{ H = "0"; J = "10625.0" }
    14 = @{ H = "0"; J = "33660.0" }
    15 = @{ H = "0"; J = "3911.7" }
    16 = @{ H = "0"; J = "816.0" }
    17 = @{ H = "0"; J = "57120.0" }
    18 = @{ H = "0"; J = "799.0" }
    19 = @{ H = "0"; J = "4896.0" }
    20 = @{ H = "0"; J = "63648.0" }
    21 = @{ H = "0"; J = "6936.0" }
    22 = @{ H = "0"; J = "8078.4" }
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Range("H$row").Value = "'" + $values.H
    $ws.Range("J$row").Value = "'" + $values.J
}

Write-Host "Updated Horas Extra / Salario Total for rows 2-22"
